$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA_FILE")

# ---------------------------------------------------------------------
# 0) Remove all existing hyperlinks up-front (the engine's Hyperlinks
#    collection Delete() call removes every hyperlink on the sheet, so
#    we do this once and re-add the two that must survive at the end).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 1) Drop columns H:J (content + widths). Deleting the columns shifts
#    the old (blank, bold-styled) column K into column H, which is
#    exactly the desired end state for H1.
# ---------------------------------------------------------------------
$ws.Range("H1:J1").EntireColumn.Delete()

# ---------------------------------------------------------------------
# 2) Wipe out the rows that disappear completely (row 7, rows 9-16).
# ---------------------------------------------------------------------
$ws.Rows.Item(7).Clear()
$ws.Rows.Item(9).Clear()
$ws.Rows.Item(10).Clear()
$ws.Rows.Item(11).Clear()
$ws.Rows.Item(12).Clear()
$ws.Rows.Item(13).Clear()
$ws.Rows.Item(14).Clear()
$ws.Rows.Item(15).Clear()
$ws.Rows.Item(16).Clear()

# row 20 disappears entirely as well
$ws.Rows.Item(20).Clear()

# ---------------------------------------------------------------------
# 3) Remove the stray cells inside the rows that survive.
# ---------------------------------------------------------------------
$ws.Range("C2").Clear()
$ws.Range("G2").Clear()
$ws.Range("H2").Clear()

$ws.Range("H3").Clear()

$ws.Range("D4").Clear()
$ws.Range("F4").Clear()

$ws.Range("E5").Clear()
$ws.Range("J5").Clear()

$ws.Range("E6").Clear()
$ws.Range("J6").Clear()

$ws.Range("B18").Clear()
$ws.Range("B19").Clear()

# ---------------------------------------------------------------------
# 4) Apply the s=6 "bold-less but explicit font" style used by column A
#    test-case names onto the rows that need it (copy the format that
#    already lives on A4 so we reuse cellXf index 6).
# ---------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5) Set the new test-case / browser text for rows 2-6.
# ---------------------------------------------------------------------
$ws.Range("A2").Value2 = "TC_001_LOGIN"
$ws.Range("B2").Value2 = "Chrome"

$ws.Range("A3").Value2 = "TC_002_VERIFY_HOMEPAGE"
$ws.Range("B3").Value2 = "Chrome"

$ws.Range("A4").Value2 = "TC_003_AMAZON_SEARCH ITEM"
$ws.Range("B4").Value2 = "Chrome"
$ws.Range("E4").Value2 = "laptop"

$ws.Range("A5").Value2 = "TC_004_AMAZON_ADD TO CART"
$ws.Range("B5").Value2 = "Chrome"

$ws.Range("A6").Value2 = "TC_005_PROCEED_TO_BUY"
$ws.Range("B6").Value2 = "Chrome"

# Row 4 keeps its 15pt custom row height (inherited from the template
# rows that used to sit at r=2/5); make sure rows 2,3,5,6 stay default.
$ws.Rows.Item(2).RowHeight = 14.5
$ws.Rows.Item(3).RowHeight = 14.5
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 14.5
$ws.Rows.Item(6).RowHeight = 14.5

# ---------------------------------------------------------------------
# 6) Set the login credential values (D2/F2) and make sure D3/F3/F8
#    stay present but blank with the hyperlink-style formatting (s=2),
#    by copying the format from the pre-existing F4 style-2 cell before
#    it got cleared is no longer possible, so copy from D2 itself once
#    it is restyled below.
# ---------------------------------------------------------------------
$ws.Range("D2").Value2 = "Amruta_18;"
$ws.Range("F2").Value2 = "amrutadanawade18@gmail.com"

# D2/F2/D3/F3/F8 must carry cellXf 2 (the hyperlink look). Grab that
# format from the still-intact D3 cell (already s=2 from the original
# workbook) before giving it new content.
$ws.Range("D3").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("D3").ClearContents()
$ws.Range("F3").ClearContents()

# ---------------------------------------------------------------------
# 7) Re-create the two surviving hyperlinks (targets identical to the
#    ones that used to live on D2/F2 before the edit).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:supriya@36", "", "", "supriya@36")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:supriya4@deloitte.com")

# Re-apply the hyperlink cell style (s=2) that Hyperlinks.Add() always
# overwrites with a brand new (duplicate) style index.
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4123)
$ws.Range("F2").PasteSpecial(-4123)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 8) Sheet view tidy-up: selection moves to D2, no frozen/top-left cell.
# ---------------------------------------------------------------------
$ws.Range("D2").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1

$wb.Save()
